# Weekly update: insert a new price observation row for
# "Feria Lagunitas de Puerto Montt - Ají" ahead of the existing row 342,
# pushing the rest of the block (old rows 342-405) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 342; existing rows 342:405 shift to 343:406
# and the sheet's used-range dimension grows to A1:R406 automatically.
$ws.Rows.Item(342).Insert()

# Populate the newly inserted row 342 with the new weekly record.
$ws.Cells.Item(342, 1).Value = 4
$ws.Cells.Item(342, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(342, 3).Value = "Los Lagos"
$ws.Cells.Item(342, 4).Value = 45034
$ws.Cells.Item(342, 5).Value = 10
$ws.Cells.Item(342, 6).Value = 100112021
$ws.Cells.Item(342, 7).Value = "Ají"
$ws.Cells.Item(342, 8).Value = "Inferno"
$ws.Cells.Item(342, 9).Value = "Primera"
$ws.Cells.Item(342, 10).Value = 180
$ws.Cells.Item(342, 11).Value = 25000
$ws.Cells.Item(342, 12).Value = 25000
$ws.Cells.Item(342, 13).Value = 25000
$ws.Cells.Item(342, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(342, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(342, 16).Value = 2500
$ws.Cells.Item(342, 17).Value = 10
$ws.Cells.Item(342, 18).Value = "Hortaliza"
